$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 7.774617366771236
$ws.Range("D2").Value = 2.807301892357337
$ws.Range("E2").Value = 10.90867421731391
$ws.Range("F2").Value = 21.55254841683206
$ws.Range("G2").Value = 3.570449054546622
$ws.Range("M2").Value = 17.20267214287284
$ws.Range("N2").Value = 17.48603232932477
$ws.Range("O2").Value = 18.50298882304369

# Row 3
$ws.Range("B3").Value = 7.700555740213737
$ws.Range("D3").Value = 2.814435072826421
$ws.Range("E3").Value = 11.04947273688814
$ws.Range("F3").Value = 21.09981515155394
$ws.Range("G3").Value = 3.57342552290726
$ws.Range("M3").Value = 16.51876177899656
$ws.Range("N3").Value = 17.41378050168197
$ws.Range("O3").Value = 18.23520480292951

# Row 4
$ws.Range("B4").Value = 7.656515638299742
$ws.Range("D4").Value = 2.818976749423928
$ws.Range("E4").Value = 11.14083218709072
$ws.Range("F4").Value = 20.82488463417977
$ws.Range("G4").Value = 3.5753496926186
$ws.Range("M4").Value = 16.08583824438596
$ws.Range("N4").Value = 17.37178291655026
$ws.Range("O4").Value = 18.07535155232397

# Row 5
$ws.Range("B5").Value = 7.63895006430598
$ws.Range("D5").Value = 2.820868485459677
$ws.Range("E5").Value = 11.17929201775843
$ws.Range("F5").Value = 20.71379586405109
$ws.Range("G5").Value = 3.576158186909787
$ws.Range("M5").Value = 15.90640376780496
$ws.Range("N5").Value = 17.35527741066324
$ws.Range("O5").Value = 18.01144519144099

# Row 6
$ws.Range("B6").Value = 7.636056919364438
$ws.Range("D6").Value = 2.821185088710054
$ws.Range("E6").Value = 11.18575242338049
$ws.Range("F6").Value = 20.69541219659289
$ws.Range("G6").Value = 3.576293911793079
$ws.Range("M6").Value = 15.87643474762161
$ws.Range("N6").Value = 17.35257385668458
$ws.Range("O6").Value = 18.00091066463538

# Row 7
$ws.Range("B7").Value = 7.656277173053823
$ws.Range("D7").Value = 2.819002095842181
$ws.Range("E7").Value = 11.1413458944926
$ws.Range("F7").Value = 20.82338237561782
$ws.Range("G7").Value = 3.575360497428602
$ws.Range("M7").Value = 16.08343017663995
$ws.Range("N7").Value = 17.37155783416786
$ws.Range("O7").Value = 18.07448457545376

# Row 8
$ws.Range("B8").Value = 7.748795800816517
$ws.Range("D8").Value = 2.80972797308185
$ws.Range("E8").Value = 10.95620014965819
$ws.Range("F8").Value = 21.395921077485
$ws.Range("G8").Value = 3.571455342767219
$ws.Range("M8").Value = 16.96970147104308
$ws.Range("N8").Value = 17.46063758714847
$ws.Range("O8").Value = 18.40976177986165

# Row 9
$ws.Range("B9").Value = 7.940670946557532
$ws.Range("D9").Value = 2.79281392146622
$ws.Range("E9").Value = 10.63227967189016
$ws.Range("F9").Value = 22.53503057811182
$ws.Range("G9").Value = 3.564559903170369
$ws.Range("M9").Value = 18.5948118151233
$ws.Range("N9").Value = 17.65344821424603
$ws.Range("O9").Value = 19.09955375347315

# Row 10
$ws.Range("B10").Value = 8.086728820589148
$ws.Range("D10").Value = 2.781146184999824
$ws.Range("E10").Value = 10.41846437706364
$ws.Range("F10").Value = 23.37171630732966
$ws.Range("G10").Value = 3.559953150715788
$ws.Range("M10").Value = 19.70884273577429
$ws.Range("N10").Value = 17.80526283336567
$ws.Range("O10").Value = 19.62069452956872

# Row 11
$ws.Range("B11").Value = 8.154001092441652
$ws.Range("D11").Value = 2.77599949369911
$ws.Range("E11").Value = 10.32651768384583
$ws.Range("F11").Value = 23.75026479306381
$ws.Range("G11").Value = 3.55795597539989
$ws.Range("M11").Value = 20.19641923412876
$ws.Range("N11").Value = 17.87633291535385
$ws.Range("O11").Value = 19.85975282105608

# Row 12
$ws.Range("B12").Value = 8.179571329865439
$ws.Range("D12").Value = 2.774073461159323
$ws.Range("E12").Value = 10.29247201273258
$ws.Range("F12").Value = 23.89315589817426
$ws.Range("G12").Value = 3.557213764828572
$ws.Range("M12").Value = 20.37816307346961
$ws.Range("N12").Value = 17.90351670524025
$ws.Range("O12").Value = 19.95047042970693

# Row 13
$ws.Range("B13").Value = 8.174060440692747
$ws.Range("D13").Value = 2.774487251793143
$ws.Range("E13").Value = 10.29976986876206
$ws.Range("F13").Value = 23.86240451930424
$ws.Range("G13").Value = 3.557372988404348
$ws.Range("M13").Value = 20.33915172224139
$ws.Range("N13").Value = 17.89765041917062
$ws.Range("O13").Value = 19.93092573850996

# Row 14
$ws.Range("B14").Value = 8.156102988388842
$ws.Range("D14").Value = 2.775840580142205
$ws.Range("E14").Value = 10.32370119673243
$ws.Range("F14").Value = 23.76203057387517
$ws.Range("G14").Value = 3.557894631644944
$ws.Range("M14").Value = 20.21143002158113
$ws.Range("N14").Value = 17.87856401195317
$ws.Range("O14").Value = 19.86721291092667

# Row 15
$ws.Range("B15").Value = 8.145115292701211
$ws.Range("D15").Value = 2.776672508818987
$ws.Range("E15").Value = 10.33846067336583
$ws.Range("F15").Value = 23.7004845229105
$ws.Range("G15").Value = 3.558215983927008
$ws.Range("M15").Value = 20.13281679586349
$ws.Range("N15").Value = 17.86690780764875
$ws.Range("O15").Value = 19.82820912296607

# Row 16
$ws.Range("B16").Value = 8.08234709381907
$ws.Range("D16").Value = 2.781485752362093
$ws.Range("E16").Value = 10.42458098345059
$ws.Range("F16").Value = 23.34692183919134
$ws.Range("G16").Value = 3.560085645058084
$ws.Range("M16").Value = 19.67658056617493
$ws.Range("N16").Value = 17.8006571472379
$ws.Range("O16").Value = 19.60510339729072

# Row 17
$ws.Range("B17").Value = 8.044036517820548
$ws.Range("D17").Value = 2.784479586079914
$ws.Range("E17").Value = 10.47878058820306
$ws.Range("F17").Value = 23.12937710077227
$ws.Range("G17").Value = 3.561257780964176
$ws.Range("M17").Value = 19.39167810747271
$ws.Range("N17").Value = 17.76051648344848
$ws.Range("O17").Value = 19.46867672847427

# Row 18
$ws.Range("B18").Value = 8.022080804529349
$ws.Range("D18").Value = 2.786216730651957
$ws.Range("E18").Value = 10.51045491533344
$ws.Range("F18").Value = 23.00406587107371
$ws.Range("G18").Value = 3.561941234882622
$ws.Range("M18").Value = 19.22600886092406
$ws.Range("N18").Value = 17.7376188039983
$ws.Range("O18").Value = 19.39040045894099

# Row 19
$ws.Range("B19").Value = 8.014661356521581
$ws.Range("D19").Value = 2.786807511077783
$ws.Range("E19").Value = 10.52126497510187
$ws.Range("F19").Value = 22.96161085399122
$ws.Range("G19").Value = 3.562174235675403
$ws.Range("M19").Value = 19.1696110566508
$ws.Range("N19").Value = 17.72989924114001
$ws.Range("O19").Value = 19.36393347629684

# Row 20
$ws.Range("B20").Value = 8.048106686015561
$ws.Range("D20").Value = 2.784159319101611
$ws.Range("E20").Value = 10.47295913083346
$ws.Range("F20").Value = 23.15255547169192
$ws.Range("G20").Value = 3.561132046041612
$ws.Range("M20").Value = 19.42219384406066
$ws.Range("N20").Value = 17.76476997169418
$ws.Range("O20").Value = 19.48318032384694

# Row 21
$ws.Range("B21").Value = 8.161375118110781
$ws.Range("D21").Value = 2.775442455042681
$ws.Range("E21").Value = 10.31665095450769
$ws.Range("F21").Value = 23.79152646832162
$ws.Range("G21").Value = 3.557741031003212
$ws.Range("M21").Value = 20.24902437753957
$ws.Range("N21").Value = 17.88416293827346
$ws.Range("O21").Value = 19.88592246060505

# Row 22
$ws.Range("B22").Value = 8.235949498577815
$ws.Range("D22").Value = 2.769878891899912
$ws.Range("E22").Value = 10.21900154218994
$ws.Range("F22").Value = 24.20639956028034
$ws.Range("G22").Value = 3.555606816823433
$ws.Range("M22").Value = 20.77250868841132
$ws.Range("N22").Value = 17.96376482283565
$ws.Range("O22").Value = 20.1502112648158

# Row 23
$ws.Range("B23").Value = 8.196105469824898
$ws.Range("D23").Value = 2.772836143948089
$ws.Range("E23").Value = 10.27070377259156
$ws.Range("F23").Value = 23.98527442445974
$ws.Range("G23").Value = 3.556738410130897
$ws.Range("M23").Value = 20.4946994765034
$ws.Range("N23").Value = 17.9211419411605
$ws.Range("O23").Value = 20.00908778324155

# Row 24
$ws.Range("B24").Value = 8.04626634455273
$ws.Range("D24").Value = 2.784304062085857
$ws.Range("E24").Value = 10.47558941041408
$ws.Range("F24").Value = 23.14207727121254
$ws.Range("G24").Value = 3.561188860960972
$ws.Range("M24").Value = 19.40840350980891
$ws.Range("N24").Value = 17.76284640850066
$ws.Range("O24").Value = 19.47662275133499

# Row 25
$ws.Range("B25").Value = 7.887776272871257
$ws.Range("D25").Value = 2.797255131504068
$ws.Range("E25").Value = 10.71568813759481
$ws.Range("F25").Value = 22.2261940628728
$ws.Range("G25").Value = 3.566344240549889
$ws.Range("M25").Value = 18.16850626290467
$ws.Range("N25").Value = 17.59943762190663
$ws.Range("O25").Value = 18.91003769845784
